$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 430, pushing the existing row 430 (and all
# rows below it) down by one. This mirrors the diff: a brand-new record is
# inserted into the daily price log, and every subsequent record keeps its
# original data but moves down one row (the former last row, 532, becomes
# row 533).
$ws.Rows(430).Insert()

# Populate the newly inserted row 430 with the new record's data.
$ws.Cells.Item(430, 1).Value  = 10
$ws.Cells.Item(430, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(430, 3).Value  = "La Araucanía"
$ws.Cells.Item(430, 4).Value  = 44754
$ws.Cells.Item(430, 5).Value  = 9
$ws.Cells.Item(430, 6).Value  = 100112043
$ws.Cells.Item(430, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(430, 8).Value  = "Sin especificar"
$ws.Cells.Item(430, 9).Value  = "Primera"
$ws.Cells.Item(430, 10).Value = 170
$ws.Cells.Item(430, 11).Value = 20000
$ws.Cells.Item(430, 12).Value = 21000
$ws.Cells.Item(430, 13).Value = 20471
$ws.Cells.Item(430, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(430, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(430, 16).Value = 409
$ws.Cells.Item(430, 17).Value = 50
$ws.Cells.Item(430, 18).Value = "Hortaliza"
